$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D/E values in this sheet are literal text (prices/percent strings, not
# real numbers), so every assignment below is prefixed with a leading apostrophe
# to stop Excel from auto-converting numeric-looking text (e.g. "1.000") into a
# number (which would also silently drop trailing/insignificant zeros).

$ws.Range("D2").Value = '''28.250.56'
$ws.Range("E2").Value = '''  +0.80%  '
$ws.Range("D3").Value = '''1.872.46'
$ws.Range("E3").Value = '''  +3.59%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '''  -0.36%  '
$ws.Range("D5").Value = '''311.50'
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '''  -0.67%  '
$ws.Range("D7").Value = '''0.5077'
$ws.Range("E7").Value = '''  +1.88%  '
$ws.Range("D8").Value = '''0.3922'
$ws.Range("E8").Value = '''  +1.56%  '
$ws.Range("D9").Value = '''0.09715'
$ws.Range("E9").Value = '''  +3.91%  '
$ws.Range("D10").Value = '''1.143'
$ws.Range("E10").Value = '''  +4.32%  '
$ws.Range("D11").Value = '''40.93'
$ws.Range("E11").Value = '''  +0.73%  '
$ws.Range("D12").Value = '''6.515'
$ws.Range("E12").Value = '''  +2.72%  '
$ws.Range("D13").Value = '''21.02'
$ws.Range("E13").Value = '''  +1.46%  '
$ws.Range("D14").Value = '''1.872.76'
$ws.Range("E14").Value = '''  +3.68%  '
$ws.Range("D15").Value = '''7.449'
$ws.Range("E15").Value = '''  +3.08%  '
$ws.Range("E16").Value = '''  -0.36%  '
$ws.Range("E17").Value = '''  +1.63%  '
$ws.Range("D18").Value = '''93.06'
$ws.Range("E18").Value = '''  +1.01%  '
$ws.Range("D19").Value = '''0.06591'
$ws.Range("E19").Value = '''  +0.16%  '
$ws.Range("D20").Value = '''17.58'
$ws.Range("E20").Value = '''  +2.68%  '
$ws.Range("E21").Value = '''  -0.62%  '
$ws.Range("D22").Value = '''6.174'
$ws.Range("E22").Value = '''  +3.49%  '
$ws.Range("D23").Value = '''28.307.29'
$ws.Range("E23").Value = '''  +0.86%  '
$ws.Range("D24").Value = '''11.34'
$ws.Range("E24").Value = '''  +3.00%  '
$ws.Range("D25").Value = '''2.287'
$ws.Range("E25").Value = '''  +2.43%  '
$ws.Range("D26").Value = '''2.554'
$ws.Range("E26").Value = '''  +7.48%  '
$ws.Range("D27").Value = '''2.089.88'
$ws.Range("E27").Value = '''  +3.67%  '
$ws.Range("D28").Value = '''21.24'
$ws.Range("E28").Value = '''  +3.72%  '
$ws.Range("D29").Value = '''158.44'
$ws.Range("E29").Value = '''  -0.14%  '
$ws.Range("D30").Value = '''127.71'
$ws.Range("E30").Value = '''  +0.39%  '
$ws.Range("D31").Value = '''0.1063'
$ws.Range("E31").Value = '''  -0.83%  '
$ws.Range("D32").Value = '''1.071'
$ws.Range("E32").Value = '''  +2.41%  '
$ws.Range("D33").Value = '''5.642'
$ws.Range("E33").Value = '''  +1.60%  '
$ws.Range("D34").Value = '''3.624'
$ws.Range("D35").Value = '''9.565'
$ws.Range("E35").Value = '''  +7.97%  '
$ws.Range("D36").Value = '''0.06726'
$ws.Range("E36").Value = '''  -1.75%  '
$ws.Range("D37").Value = '''0.02386'
$ws.Range("E37").Value = '''  +3.84%  '
$ws.Range("D38").Value = '''0.2191'
$ws.Range("E38").Value = '''  +2.61%  '

# Rows 39 and 40 swap ranking positions (Aptos moves above TheSandbox),
# with refreshed price/volume figures.
$ws.Range("B39").Value = '''Aptos'
$ws.Range("C39").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''11.52'
$ws.Range("E39").Value = '''  +1.35%  '

$ws.Range("B40").Value = '''TheSandbox'
$ws.Range("C40").Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6379'
$ws.Range("E40").Value = '''  +3.95%  '

$ws.Range("D41").Value = '''4.981'
$ws.Range("E41").Value = '''  +0.94%  '
$ws.Range("D42").Value = '''1.186'
$ws.Range("E42").Value = '''  +3.57%  '
$ws.Range("E43").Value = '''  -0.61%  '
$ws.Range("E44").Value = '''  +3.56%  '
$ws.Range("D45").Value = '''0.6020'
$ws.Range("E45").Value = '''  +2.42%  '
$ws.Range("D46").Value = '''3.662'
$ws.Range("E46").Value = '''  +0.00%  '
$ws.Range("E47").Value = '''  -1.97%  '
$ws.Range("D48").Value = '''2.001'
$ws.Range("E48").Value = '''  +2.84%  '
$ws.Range("D49").Value = '''124.05'
$ws.Range("E49").Value = '''  -0.01%  '
$ws.Range("E50").Value = '''  +2.12%  '
$ws.Range("E51").Value = '''  +1.42%  '
